$d = $word.ActiveDocument

# The presentation date was corrected from 17/03/2022 to 28/04/2022.
# (In the authored edit this was typed over in two selections -
#  "17" -> "28" and the "3" in "03" -> "4" - which is replicated below
#  so the resulting text matches exactly: "Apresentação 28/04/2022;".)

$rng = $d.Content
$found = $rng.Find.Execute("17/03/2022;", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Replace "17" with "28"
    $day = $d.Range($rng.Start, $rng.Start + 2)
    $day.Text = "28"

    # Replace the "3" in "03" with "4" -> "04"
    $month = $d.Range($rng.Start + 4, $rng.Start + 5)
    $month.Text = "4"
} else {
    # Fallback: direct replace in case the run text ever changes shape.
    $d.Content.Find.Execute("Apresentação 17/03/2022;", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "Apresentação 28/04/2022;", 2)
}
